# BOM.xlsx update — "aanpassingen bom, gerber files voor bestelling"
#
# 1. G14: mark the zener-diode "vraag aan stijn" note as resolved.
# 2. G29: the zener diode is now sourced from Comchip (MMSZ4678-HF) instead
#    of ON Semiconductor (MMSZ4678T1G); add the corresponding hyperlink.
# 3. G42: TPS563200 buck converter now ordered from sinuss.be instead of
#    mouser.
# 4. Add a new "Bestelling" (order) summary block below the BOM table,
#    listing the 4 orders placed (sinuss, pcb, mouser, aliexpress) with
#    their cost, and a totaal (total) row that sums them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. G14 ---------------------------------------------------------------
$ws.Range("G14").Value = "aan stijn vragen (in orde)"

# --- 2. G29 -----------------------------------------------------------------
$ws.Range("G29").Value = "https://www.mouser.be/ProductDetail/Comchip-Technology/MMSZ4678-HF?qs=GBLSl2AkirsCXt3Ql7prKg%3D%3D"
$ws.Hyperlinks.Add($ws.Range("G29"), "https://www.mouser.be/ProductDetail/Comchip-Technology/MMSZ4678-HF?qs=GBLSl2AkirsCXt3Ql7prKg%3D%3D") | Out-Null

# --- 3. G42 (text only – keep the existing hyperlink relationship as-is) ---
$ws.Range("G42").Value = "https://sinuss.be/componenten/halfgeleiders-sensors/ics/voltage-regulators/dc-dc-switching-regulators-adjustable/2450170-tps563200ddct-dcdc-conv-sync-buck-650khz-sot-23-6-texas-instruments"

# --- 4. New "Bestelling" block --------------------------------------------
$ws.Range("A46").Value = "Bestelling"

$ws.Range("A47").Value = "sinuss"
$ws.Range("B47").Value = 10.49
$ws.Range("D47").Value = "besteld"

$ws.Range("A48").Value = "pcb"
$ws.Range("B48").Value = 56.85
$ws.Range("D48").Value = "besteld"

$ws.Range("A49").Value = "mouser"
$ws.Range("B49").Value = 139.82
$ws.Range("D49").Value = "besteld"

$ws.Range("A50").Value = "aliexpress"
$ws.Range("B50").Value = 11.3
$ws.Range("D50").Value = "besteld"

$ws.Range("A52").Value = "totaal"
$ws.Range("B52").Formula = "=B47+B48+B49+B50"

# --- view state: selection on B53, scrolled near the new block ------------
$ws.Range("B53").Select()
